# "errors de canvi de dni" - fix the player data that got mixed up when the
# dni values were corrected: names need proper capitalisation, the
# tarjeta_sanitaria numbers are really text (they contain a slash), the
# mobil/telefon numbers were placeholders, numero_soci was wrong, the
# correu_electronic column needs to hold real e-mail addresses (as
# clickable mailto hyperlinks), and a missing indata_naixement (F5) needs
# to be filled in.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- nom (B2:B9): capitalise properly ---------------------------------
$ws.Range("B2").Value = "Joan"
$ws.Range("B3").Value = "Josep"
$ws.Range("B4").Value = "Pere"
$ws.Range("B5").Value = "Antoni"
$ws.Range("B6").Value = "Sergi"
$ws.Range("B7").Value = "David"
$ws.Range("B8").Value = "Joel"
$ws.Range("B9").Value = "Miki"

# --- correu_electronic (K2:K9): real e-mail addresses as hyperlinks ---
$ws.Range("K2").Value = "joan.noguer@hotmail.com"
$ws.Hyperlinks.Add($ws.Range("K2"), "mailto:joan.noguer@hotmail.com")

$ws.Range("K3").Value = "josep1@gmail.com"
$ws.Hyperlinks.Add($ws.Range("K3"), "mailto:josep1@gmail.com")

$ws.Range("K4").Value = "pereee@hotmail.com"
$ws.Hyperlinks.Add($ws.Range("K4"), "mailto:pereee@hotmail.com")

$ws.Range("K5").Value = "antoniiooooo@hotmail.com"
$ws.Hyperlinks.Add($ws.Range("K5"), "mailto:antoniiooooo@hotmail.com")

$ws.Range("K6").Value = "sergi92@gmail.com"
$ws.Hyperlinks.Add($ws.Range("K6"), "mailto:sergi92@gmail.com")

$ws.Range("K7").Value = "sisterna@hotmail.com"
$ws.Hyperlinks.Add($ws.Range("K7"), "mailto:sisterna@hotmail.com")

$ws.Range("K8").Value = "joey@yahoo.es"
$ws.Hyperlinks.Add($ws.Range("K8"), "mailto:joey@yahoo.es")

$ws.Range("K9").Value = "yuyuhakusho@niancat.japo"
$ws.Hyperlinks.Add($ws.Range("K9"), "mailto:yuyuhakusho@niancat.japo")

# --- tarjeta_sanitaria (G2:G9): these are text (contain "/"), not numbers
$ws.Range("G2").Value = "39485615/76"
$ws.Range("G3").Value = "123648679/13"
$ws.Range("G4").Value = "75492367/09"
$ws.Range("G5").Value = "75892345/68"
$ws.Range("G6").Value = "909873256/89"
$ws.Range("G7").Value = "37584979/99"
$ws.Range("G8").Value = "647836578/78"
$ws.Range("G9").Value = "176947899/13"

# --- indata_naixement (F2:F9): corrected birth dates (F5 was missing) -
$ws.Range("F2").Value = 36679
$ws.Range("F3").Value = 36408
$ws.Range("F4").Value = 37667
$ws.Range("F5").Value = 34202
$ws.Range("F6").Value = 37429
$ws.Range("F7").Value = 35492
$ws.Range("F8").Value = 35286
$ws.Range("F9").Value = 38838

# --- mobil (I6:I9): replace placeholder numbers ------------------------
$ws.Range("I6").Value = 68743598
$ws.Range("I7").Value = 623479043
$ws.Range("I8").Value = 612879544
$ws.Range("I9").Value = 680754376

# --- telefon (J2:J9): replace placeholder numbers -----------------------
$ws.Range("J2").Value = 972867698
$ws.Range("J3").Value = 972866987
$ws.Range("J4").Value = 972860970
$ws.Range("J5").Value = 972863687
$ws.Range("J6").Value = 972868843
$ws.Range("J7").Value = 972867735
$ws.Range("J8").Value = 972860354
$ws.Range("J9").Value = 972864412

# --- numero_soci (L2:L9): corrected membership numbers ------------------
$ws.Range("L2").Value = 15
$ws.Range("L3").Value = 166
$ws.Range("L4").Value = 17
$ws.Range("L5").Value = 655
$ws.Range("L6").Value = 78
$ws.Range("L7").Value = 89
$ws.Range("L8").Value = 56
$ws.Range("L9").Value = 88

# --- restore the selection to where the author left it ------------------
$ws.Range("F10").Select() | Out-Null
